# Applies the price/volume/coin updates from the Oct 23 2023 GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.832.16"
$ws.Range("E2").Value = "  +3.16%  "

# Row 3
$ws.Range("D3").Value = "1.679.68"
$ws.Range("E3").Value = "  +3.15%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.20%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.41%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.536"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.38%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.13%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "28.97"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.65%  "

# Row 9
$ws.Range("E9").Value = "  +2.42%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0645"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.17%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0906"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.59%  "

# Row 12
$ws.Range("D12").Value = "1.921.56"
$ws.Range("E12").Value = "  +3.25%  "

# Row 13
$ws.Range("D13").Value = "1.685.04"
$ws.Range("E13").Value = "  +3.45%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +10.09%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.604"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.55%  "

# Row 16
$ws.Range("E16").Value = "  +5.54%  "

# Row 17
$ws.Range("D17").Value = "30.799.47"
$ws.Range("E17").Value = "  +2.97%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "66.05"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.19%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "243.46"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.21%  "

# Row 20
$ws.Range("E20").Value = "  +2.70%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.998"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.26%  "

# Row 22
$ws.Range("E22").Value = "  +2.97%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.96"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.97%  "

# Row 24
$ws.Range("E24").Value = "  +0.26%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.23"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.90%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.84"
$ws.Range("D26").Style = "Normal"

# Row 27
$ws.Range("E27").Value = "  +2.33%  "

# Row 28
$ws.Range("E28").Value = "  +2.35%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.15%  "

# Row 30
$ws.Range("E30").Value = "  +1.41%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.14"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.85%  "

# Row 32
$ws.Range("E32").Value = "  +3.07%  "

# Row 33
$ws.Range("D33").Value = "1.523.42"
$ws.Range("E33").Value = "  +6.90%  "

# Row 34
$ws.Range("E34").Value = "  +4.32%  "

# Row 35
$ws.Range("E35").Value = "  +5.36%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "84.03"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +12.45%  "

# Row 37
$ws.Range("E37").Value = "  +0.73%  "

# Row 38
$ws.Range("E38").Value = "  +8.55%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0178"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.03%  "

# Row 40
$ws.Range("B40").Value = "HuobiToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.29"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.10%  "

# Row 41
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.65"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.75%  "

# Row 42
$ws.Range("E42").Value = "  +2.56%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.839"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.57%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0500"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.58%  "

# Row 45
$ws.Range("E45").Value = "  +2.00%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.08%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.54"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.93%  "

# Row 48
$ws.Range("B48").Value = "BitcoinSV"
$ws.Range("C48").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "50.71"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.06%  "

# Row 49
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "1.812.16"
$ws.Range("E49").Value = "  +2.44%  "

# Row 50
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0117"
$ws.Range("E50").Value = "  +5.17%  "

# Row 51
$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "92.84"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.24%  "
